$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in N3:N111 while preserving the cell formatting/style.
$ws.Range("N3:N111").ClearContents()
